$wb = $excel.ActiveWorkbook
$changes = @(
  @{Sheet="ALC"; Cell="H40"; Value=1474.5},
  @{Sheet="ALC"; Cell="I40"; Value=0},
  @{Sheet="ALC"; Cell="J40"; Value=1474.5},
  @{Sheet="ALC"; Cell="K40"; Value=0},
  @{Sheet="ALC"; Cell="L40"; Value=1474.5},
  @{Sheet="ALC"; Cell="M40"; Value=$null},
  @{Sheet="ALC"; Cell="N40"; Value=-1824.5},
  @{Sheet="ALC"; Cell="H76"; Value=4666.3335},
  @{Sheet="ALC"; Cell="J76"; Value=3999},
  @{Sheet="ALC"; Cell="L76"; Value=3999},
  @{Sheet="ALC"; Cell="N76"; Value=-4629},
  @{Sheet="ALC"; Cell="H79"; Value=4666.3335},
  @{Sheet="ALC"; Cell="J79"; Value=3999},
  @{Sheet="ALC"; Cell="L79"; Value=3999},
  @{Sheet="ALC"; Cell="N79"; Value=-6183},
  @{Sheet="ALC"; Cell="H100"; Value=6064.8335},
  @{Sheet="ALC"; Cell="I100"; Value=4131.6665},
  @{Sheet="ALC"; Cell="J100"; Value=7998},
  @{Sheet="ALC"; Cell="K100"; Value=4131.6665},
  @{Sheet="ALC"; Cell="L100"; Value=7998},
  @{Sheet="ALC"; Cell="M100"; Value=-3590.6665},
  @{Sheet="ALC"; Cell="N100"; Value=-9080},
  @{Sheet="ALC"; Cell="H112"; Value=1578.8966},
  @{Sheet="ALC"; Cell="J112"; Value=1655.7916},
  @{Sheet="ALC"; Cell="L112"; Value=4967.3748},
  @{Sheet="ALC"; Cell="N112"; Value=-7183.3748},
  @{Sheet="ALC"; Cell="H118"; Value=810.5},
  @{Sheet="ALC"; Cell="I118"; Value=786.38464},
  @{Sheet="ALC"; Cell="J118"; Value=915},
  @{Sheet="ALC"; Cell="K118"; Value=2359.15392},
  @{Sheet="ALC"; Cell="L118"; Value=2745},
  @{Sheet="ALC"; Cell="M118"; Value=-702.1539199999997},
  @{Sheet="ALC"; Cell="N118"; Value=-6059},
  @{Sheet="ALC"; Cell="H137"; Value=1482.6154},
  @{Sheet="ALC"; Cell="I137"; Value=1482.6154},
  @{Sheet="ALC"; Cell="J137"; Value=0},
  @{Sheet="ALC"; Cell="K137"; Value=4447.8462},
  @{Sheet="ALC"; Cell="L137"; Value=0},
  @{Sheet="ALC"; Cell="M137"; Value=-1897.8462},
  @{Sheet="ALC"; Cell="N137"; Value=$null},
  @{Sheet="ARM"; Cell="H61"; Value=7337},
  @{Sheet="ARM"; Cell="I61"; Value=7337},
  @{Sheet="ARM"; Cell="K61"; Value=7337},
  @{Sheet="ARM"; Cell="M61"; Value=-7125},
  @{Sheet="ARM"; Cell="H102"; Value=4609.6},
  @{Sheet="ARM"; Cell="I102"; Value=4512},
  @{Sheet="ARM"; Cell="K102"; Value=4512},
  @{Sheet="ARM"; Cell="M102"; Value=-2890},
  @{Sheet="ARM"; Cell="H106"; Value=26184.75},
  @{Sheet="ARM"; Cell="J106"; Value=26184.75},
  @{Sheet="ARM"; Cell="L106"; Value=26184.75},
  @{Sheet="ARM"; Cell="N106"; Value=-28708.75},
  @{Sheet="ARM"; Cell="H132"; Value=1481.7},
  @{Sheet="ARM"; Cell="I132"; Value=1524.2222},
  @{Sheet="ARM"; Cell="K132"; Value=4572.6666},
  @{Sheet="ARM"; Cell="M132"; Value=-2042.6666},
  @{Sheet="ARM"; Cell="H136"; Value=7337},
  @{Sheet="ARM"; Cell="I136"; Value=7337},
  @{Sheet="ARM"; Cell="K136"; Value=22011},
  @{Sheet="ARM"; Cell="M136"; Value=-19461},
  @{Sheet="BSM"; Cell="H86"; Value=2522.3},
  @{Sheet="BSM"; Cell="I86"; Value=1824},
  @{Sheet="BSM"; Cell="K86"; Value=1824},
  @{Sheet="BSM"; Cell="M86"; Value=-701},
  @{Sheet="BSM"; Cell="H89"; Value=2522.3},
  @{Sheet="BSM"; Cell="I89"; Value=1824},
  @{Sheet="BSM"; Cell="K89"; Value=9120},
  @{Sheet="BSM"; Cell="M89"; Value=-3504},
  @{Sheet="BSM"; Cell="H94"; Value=3184.7144},
  @{Sheet="BSM"; Cell="I94"; Value=3184.7144},
  @{Sheet="BSM"; Cell="K94"; Value=3184.7144},
  @{Sheet="BSM"; Cell="M94"; Value=-2733.7144},
  @{Sheet="BSM"; Cell="H107"; Value=964.6},
  @{Sheet="BSM"; Cell="I107"; Value=999.5},
  @{Sheet="BSM"; Cell="J107"; Value=825},
  @{Sheet="BSM"; Cell="K107"; Value=999.5},
  @{Sheet="BSM"; Cell="L107"; Value=825},
  @{Sheet="BSM"; Cell="M107"; Value=920.5},
  @{Sheet="BSM"; Cell="N107"; Value=-4665},
  @{Sheet="BSM"; Cell="H134"; Value=9182.058999999999},
  @{Sheet="BSM"; Cell="I134"; Value=8737},
  @{Sheet="BSM"; Cell="J134"; Value=9998},
  @{Sheet="BSM"; Cell="K134"; Value=26211},
  @{Sheet="BSM"; Cell="L134"; Value=29994},
  @{Sheet="BSM"; Cell="M134"; Value=-23676},
  @{Sheet="BSM"; Cell="N134"; Value=-35064},
  @{Sheet="CRP"; Cell="H8"; Value=910},
  @{Sheet="CRP"; Cell="J8"; Value=910},
  @{Sheet="CRP"; Cell="L8"; Value=910},
  @{Sheet="CRP"; Cell="N8"; Value=-1190},
  @{Sheet="CRP"; Cell="H60"; Value=23333.334},
  @{Sheet="CRP"; Cell="I60"; Value=15000},
  @{Sheet="CRP"; Cell="K60"; Value=15000},
  @{Sheet="CRP"; Cell="M60"; Value=-14489},
  @{Sheet="CRP"; Cell="H97"; Value=0},
  @{Sheet="CRP"; Cell="J97"; Value=0},
  @{Sheet="CRP"; Cell="L97"; Value=0},
  @{Sheet="CRP"; Cell="N97"; Value=$null},
  @{Sheet="CRP"; Cell="H105"; Value=737.25},
  @{Sheet="CRP"; Cell="I105"; Value=744.8},
  @{Sheet="CRP"; Cell="J105"; Value=699.5},
  @{Sheet="CRP"; Cell="K105"; Value=744.8},
  @{Sheet="CRP"; Cell="L105"; Value=699.5},
  @{Sheet="CRP"; Cell="M105"; Value=1002.2},
  @{Sheet="CRP"; Cell="N105"; Value=-4193.5},
  @{Sheet="CRP"; Cell="H109"; Value=41275},
  @{Sheet="CRP"; Cell="J109"; Value=41275},
  @{Sheet="CRP"; Cell="L109"; Value=41275},
  @{Sheet="CRP"; Cell="N109"; Value=-43355},
  @{Sheet="CRP"; Cell="H132"; Value=881},
  @{Sheet="CRP"; Cell="I132"; Value=881},
  @{Sheet="CRP"; Cell="K132"; Value=2643},
  @{Sheet="CRP"; Cell="M132"; Value=-113},
  @{Sheet="CUL"; Cell="H64"; Value=1843.3334},
  @{Sheet="CUL"; Cell="I64"; Value=265},
  @{Sheet="CUL"; Cell="K64"; Value=795},
  @{Sheet="CUL"; Cell="M64"; Value=-525},
  @{Sheet="CUL"; Cell="H67"; Value=1843.3334},
  @{Sheet="CUL"; Cell="I67"; Value=265},
  @{Sheet="CUL"; Cell="K67"; Value=795},
  @{Sheet="CUL"; Cell="M67"; Value=141},
  @{Sheet="CUL"; Cell="H113"; Value=555.125},
  @{Sheet="CUL"; Cell="J113"; Value=691.5},
  @{Sheet="CUL"; Cell="L113"; Value=2074.5},
  @{Sheet="CUL"; Cell="N113"; Value=-6414.5},
  @{Sheet="CUL"; Cell="H131"; Value=992.2},
  @{Sheet="CUL"; Cell="J131"; Value=1033},
  @{Sheet="CUL"; Cell="L131"; Value=3099},
  @{Sheet="CUL"; Cell="N131"; Value=-13179},
  @{Sheet="CUL"; Cell="H137"; Value=2664},
  @{Sheet="CUL"; Cell="I137"; Value=2060.6667},
  @{Sheet="CUL"; Cell="J137"; Value=3116.5},
  @{Sheet="CUL"; Cell="K137"; Value=6182.000100000001},
  @{Sheet="CUL"; Cell="L137"; Value=9349.5},
  @{Sheet="CUL"; Cell="M137"; Value=-1082.000100000001},
  @{Sheet="CUL"; Cell="N137"; Value=-19549.5},
  @{Sheet="CUL"; Cell="H140"; Value=2412.7693},
  @{Sheet="CUL"; Cell="I140"; Value=2412.7693},
  @{Sheet="CUL"; Cell="J140"; Value=0},
  @{Sheet="CUL"; Cell="K140"; Value=7238.3079},
  @{Sheet="CUL"; Cell="L140"; Value=0},
  @{Sheet="CUL"; Cell="M140"; Value=-2058.3079},
  @{Sheet="CUL"; Cell="N140"; Value=$null},
  @{Sheet="CUL"; Cell="H141"; Value=1782.2},
  @{Sheet="CUL"; Cell="I141"; Value=1782.2},
  @{Sheet="CUL"; Cell="K141"; Value=5346.6},
  @{Sheet="CUL"; Cell="M141"; Value=-166.6000000000004},
  @{Sheet="GSM"; Cell="H2"; Value=296.3158},
  @{Sheet="GSM"; Cell="I2"; Value=326.46667},
  @{Sheet="GSM"; Cell="J2"; Value=183.25},
  @{Sheet="GSM"; Cell="K2"; Value=326.46667},
  @{Sheet="GSM"; Cell="L2"; Value=183.25},
  @{Sheet="GSM"; Cell="M2"; Value=-213.46667},
  @{Sheet="GSM"; Cell="N2"; Value=-409.25},
  @{Sheet="GSM"; Cell="H14"; Value=20003756},
  @{Sheet="GSM"; Cell="I14"; Value=20003756},
  @{Sheet="GSM"; Cell="K14"; Value=20003756},
  @{Sheet="GSM"; Cell="M14"; Value=-20003588},
  @{Sheet="GSM"; Cell="H80"; Value=2762.3845},
  @{Sheet="GSM"; Cell="I80"; Value=1633},
  @{Sheet="GSM"; Cell="J80"; Value=3264.3333},
  @{Sheet="GSM"; Cell="K80"; Value=1633},
  @{Sheet="GSM"; Cell="L80"; Value=3264.3333},
  @{Sheet="GSM"; Cell="M80"; Value=-635},
  @{Sheet="GSM"; Cell="N80"; Value=-5260.3333},
  @{Sheet="GSM"; Cell="H83"; Value=2762.3845},
  @{Sheet="GSM"; Cell="I83"; Value=1633},
  @{Sheet="GSM"; Cell="J83"; Value=3264.3333},
  @{Sheet="GSM"; Cell="K83"; Value=8165},
  @{Sheet="GSM"; Cell="L83"; Value=16321.6665},
  @{Sheet="GSM"; Cell="M83"; Value=-3173},
  @{Sheet="GSM"; Cell="N83"; Value=-26305.6665},
  @{Sheet="GSM"; Cell="H127"; Value=0},
  @{Sheet="GSM"; Cell="J127"; Value=0},
  @{Sheet="GSM"; Cell="L127"; Value=0},
  @{Sheet="GSM"; Cell="N127"; Value=$null},
  @{Sheet="GSM"; Cell="H134"; Value=96775},
  @{Sheet="GSM"; Cell="J134"; Value=96775},
  @{Sheet="GSM"; Cell="L134"; Value=290325},
  @{Sheet="GSM"; Cell="N134"; Value=-295395},
  @{Sheet="LTW"; Cell="H46"; Value=3307.4614},
  @{Sheet="LTW"; Cell="I46"; Value=2777.5557},
  @{Sheet="LTW"; Cell="J46"; Value=4499.75},
  @{Sheet="LTW"; Cell="K46"; Value=2777.5557},
  @{Sheet="LTW"; Cell="L46"; Value=4499.75},
  @{Sheet="LTW"; Cell="M46"; Value=-2589.5557},
  @{Sheet="LTW"; Cell="N46"; Value=-4875.75},
  @{Sheet="LTW"; Cell="H123"; Value=0},
  @{Sheet="LTW"; Cell="J123"; Value=0},
  @{Sheet="LTW"; Cell="L123"; Value=0},
  @{Sheet="LTW"; Cell="N123"; Value=$null},
  @{Sheet="WVR"; Cell="H41"; Value=15319.25},
  @{Sheet="WVR"; Cell="I41"; Value=16397.25},
  @{Sheet="WVR"; Cell="J41"; Value=14241.25},
  @{Sheet="WVR"; Cell="K41"; Value=16397.25},
  @{Sheet="WVR"; Cell="L41"; Value=14241.25},
  @{Sheet="WVR"; Cell="M41"; Value=-16007.25},
  @{Sheet="WVR"; Cell="N41"; Value=-15021.25},
  @{Sheet="WVR"; Cell="H81"; Value=996.6667},
  @{Sheet="WVR"; Cell="I81"; Value=996.6667},
  @{Sheet="WVR"; Cell="J81"; Value=0},
  @{Sheet="WVR"; Cell="K81"; Value=1993.3334},
  @{Sheet="WVR"; Cell="L81"; Value=0},
  @{Sheet="WVR"; Cell="M81"; Value=-932.3334},
  @{Sheet="WVR"; Cell="N81"; Value=$null},
  @{Sheet="WVR"; Cell="H84"; Value=996.6667},
  @{Sheet="WVR"; Cell="I84"; Value=996.6667},
  @{Sheet="WVR"; Cell="J84"; Value=0},
  @{Sheet="WVR"; Cell="K84"; Value=9966.666999999999},
  @{Sheet="WVR"; Cell="L84"; Value=0},
  @{Sheet="WVR"; Cell="M84"; Value=-4662.666999999999},
  @{Sheet="WVR"; Cell="N84"; Value=$null},
  @{Sheet="WVR"; Cell="H109"; Value=66500},
  @{Sheet="WVR"; Cell="J109"; Value=66500},
  @{Sheet="WVR"; Cell="L109"; Value=66500},
  @{Sheet="WVR"; Cell="N109"; Value=-69274},
  @{Sheet="WVR"; Cell="H132"; Value=1987},
  @{Sheet="WVR"; Cell="I132"; Value=1987},
  @{Sheet="WVR"; Cell="K132"; Value=5961},
  @{Sheet="WVR"; Cell="M132"; Value=-3431},
  @{Sheet="WVR"; Cell="H136"; Value=3018.4285},
  @{Sheet="WVR"; Cell="I136"; Value=3018.4285},
  @{Sheet="WVR"; Cell="K136"; Value=9055.2855},
  @{Sheet="WVR"; Cell="M136"; Value=-6505.2855}
)
foreach ($item in $changes) {
  $ws = $wb.Worksheets.Item($item.Sheet)
  if ($null -eq $item.Value) {
    $ws.Range($item.Cell).ClearContents()
  } else {
    $ws.Range($item.Cell).Value = $item.Value
  }
}